# Apply the commit's changes to debug_book.xlsx:
#  - "Compliance" sheet: add a "Notable fixes" note for the last row (F20)
#  - "Worksheet" sheet: update the DAA debug values (C/D/E columns for rows 2-4)
#    and add "Wanted"/"We got" rows (6 & 7) with updated A: 0A -> .. text
#  - Move the active selection on the "Worksheet" sheet from A7 to A6

$wb = $excel.ActiveWorkbook

$wsCompliance = $wb.Worksheets.Item("Compliance")
$wsWorksheet  = $wb.Worksheets.Item("Worksheet")

# --- Compliance sheet: new "Notable fixes" entry on row 20 ---
$wsCompliance.Range("F20").Value = "Remaining DAA problems, Z180 delivering odd results"

# --- Worksheet sheet: updated DAA test data ---
# (D column values are numeric-looking text like the original "01"/"55"/"94" --
#  prefix with an apostrophe so Excel keeps them as text, not numbers.)

# Row 2 (3F / Source)
$wsWorksheet.Range("C2").Value = "4121FA09601D59A55B8D7990020A9D29"
$wsWorksheet.Range("D2").Value = "'02"
$wsWorksheet.Range("E2").Value = "_ _ _ _ _ _ N _"

# Row 3 (CCF / SC131)
$wsWorksheet.Range("C3").Value = "4121FA09601D59A55B8D7990060A9D29"
$wsWorksheet.Range("D3").Value = "'06"
$wsWorksheet.Range("E3").Value = "_ _ _ _ _ P N _"

# Row 4 (BOX80)
$wsWorksheet.Range("C4").Value = "4121FA09601D59A55B8D799002049D29"
$wsWorksheet.Range("D4").Value = "'02"
$wsWorksheet.Range("E4").Value = "_ _ _ _ _ _ N _"

# Rows 6-7: "Wanted" / "We got" summary (replacing the old single "A: 9A -> 00" row)
$wsWorksheet.Range("B6").Value = "Wanted"
$wsWorksheet.Range("C6").Value = "A: 0A -> 0A"
$wsWorksheet.Range("B7").Value = "We got"
$wsWorksheet.Range("C7").Value = "A: 0A -> 04"

# Move the selection/active cell on the Worksheet sheet to A6 (matches the
# diff's <selection activeCell="A6" sqref="A6"/>), then re-select the
# Compliance sheet so it remains the active tab, as it was originally.
$wsWorksheet.Range("A6").Select()
$wsCompliance.Select()
